$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: copy formatting from row 34 (A=bold key style, B=value style), then set values
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("B34").Copy()
$ws.Range("B35").PasteSpecial(-4122)

# Row 36: same formatting source
$ws.Range("A34").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("B34").Copy()
$ws.Range("B36").PasteSpecial(-4122)

$ws.Range("A35").Value = "stringStart"
$ws.Range("B35").Value = "....+....1....+....2....+....3....+....4....+....5....+....6....+....7.."
$ws.Range("A36").Value = "stringEnd"
$ws.Range("B36").Value = "****** ********  End of report  ********"

$excel.CutCopyMode = $false

$ws.Range("A38").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
